$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Executive Secretary - update months/weeks/hours
$ws.Range("B2").Value = 4
$ws.Range("C2").Value = 16
$ws.Range("D2").Value = 160

# Row 3: High Council - update months/weeks, hours becomes a formula
$ws.Range("B3").Value = 5
$ws.Range("C3").Value = 20
$ws.Range("D3").Formula = "=C3*5"

# Insert a new row at position 4 for "Stake Executive Secretary"
$ws.Rows.Item(4).Insert()

$ws.Range("A4").Value = "Stake Executive Secretary"
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 4
$ws.Range("D4").Formula = "=C4*12"

# Row 5 (was row 4): Presidential Intern - now has B/C values and D is plain 471
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 4
$ws.Range("D5").Value = 471

# Row 6 (was row 5): Honors Ambassador - now has B/C values and D is plain 12
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = 4
$ws.Range("D6").Value = 12

# Row 7 (was row 6): Pre-Medical Association - D7 stays 20 (already correct after shift)
# Row 8 (was row 7): Station1 - D8 stays 100 (already correct after shift)

# Column A widened to fit the new, longer "Stake Executive Secretary" label
$ws.Columns.Item(1).ColumnWidth = 23.1

$ws.Range("G8").Select()
